# Updated master titration file
# Add a new "CTNa" column (M), shorten the "Lueker, Waters etc" label to
# "Lueker" for the first three data rows, and fill in CTNa values (14.99)
# for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column M
$ws.Range("M1").Value = "CTNa"

# Shorten carbonate_K label for rows 2-4
$ws.Range("L2").Value = "Lueker"
$ws.Range("L3").Value = "Lueker"
$ws.Range("L4").Value = "Lueker"

# New CTNa values for rows 2-8
$ws.Range("M2").Value = 14.99
$ws.Range("M3").Value = 14.99
$ws.Range("M4").Value = 14.99
$ws.Range("M5").Value = 14.99
$ws.Range("M6").Value = 14.99
$ws.Range("M7").Value = 14.99
$ws.Range("M8").Value = 14.99

# Match the recorded selection state from the edit
$ws.Range("N21").Select()
